$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source export gained a duplicate/erroneous record for account 004404248
# (PAULO, balance 1108.48). Remove that entire row so the rows below it
# (starting with account 004211368 / ILTON) shift up to fill the gap.
$target = $ws.Range("A1:A1000").Find("004404248")
if ($target -ne $null) {
    $target.EntireRow.Delete()
}
